$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet lists Kafka-Streams operators (col A) with their type-signature
# (col B) and a Status column (C). Several operators apply to both KStream
# and KTable and used to show that as two lines in a single wrapped B cell
# (e.g. "Kstream -> Kstream\nKtable -> Ktable"). This change splits each of
# those combined cells into two separate rows, merging col A down the pair
# and drawing a top/bottom divider border so the pair still reads as one
# logical entry.
# ---------------------------------------------------------------------------

# Step 1: insert a blank row under each operator that needs to be split.
# Done bottom-up so the row numbers used below stay valid while inserting.
$ws.Rows("12").Insert()   # split for old row 11 "MapValues"
$ws.Rows("10").Insert()   # split for old row 9  "GroupBy"
$ws.Rows("8").Insert()    # split for old row 7  "Foreach"
$ws.Rows("5").Insert()    # split for old row 4  "Inverse Filter"
$ws.Rows("4").Insert()    # split for old row 3  "Filter"

# Step 2: put the second line of each combined signature onto its own row,
# and trim the first line down to a single line in place.
$ws.Range("B3").Value = "Kstream -> Kstream"
$ws.Range("B4").Value = "Ktable -> Ktable"

$ws.Range("B5").Value = "Kstream -> Kstream"
$ws.Range("B6").Value = "Ktable -> Ktable"

$ws.Range("B9").Value = "KStream " + [char]8594 + " void"
$ws.Range("B10").Value = "KTable " + [char]8594 + " void"

$ws.Range("B12").Value = "KStream " + [char]8594 + " KGroupedStream"
$ws.Range("B13").Value = "KTable " + [char]8594 + " KGroupedTable"

$ws.Range("B15").Value = "KStream " + [char]8594 + " Kstream"
$ws.Range("B16").Value = "KTable " + [char]8594 + " Ktable"

# Step 3: the trimmed single-line cells no longer need to wrap (except B5,
# which keeps its original wrapped formatting - matching the source edit).
$ws.Range("B3").WrapText = $false
$ws.Range("B4").WrapText = $false
$ws.Range("B6").WrapText = $false
$ws.Range("B9").WrapText = $false
$ws.Range("B10").WrapText = $false
$ws.Range("B12").WrapText = $false
$ws.Range("B13").WrapText = $false
$ws.Range("B15").WrapText = $false
$ws.Range("B16").WrapText = $false

# Step 4: the new blank rows' C cells should carry the same plain bordered
# style as every other Status cell - copy it over from the row above.
$ws.Range("C3").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("C5").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("C9").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("C12").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("C16").PasteSpecial(-4122)

# Step 5: merge column A across each split pair, vertically-center the
# label, and draw a divider (top border on the first row of the pair,
# bottom border on the second) so the merged look is preserved visually
# even though the cell itself is only "half" of the pair.
$pairs = @(3, 5, 9, 12, 15)
foreach ($top in $pairs) {
    $bottom = $top + 1
    $ws.Range("A" + $top + ":A" + $bottom).Merge()

    $ws.Range("A" + $top).Borders.Item(9).LineStyle = 0
    $ws.Range("A" + $top).VerticalAlignment = -4108

    $ws.Range("A" + $bottom).Borders.Item(8).LineStyle = 0
    $ws.Range("A" + $bottom).VerticalAlignment = -4108
}

# Step 6: dimension / selection bookkeeping to match the final used range.
$ws.Range("C13").Select()
